$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 34, pushing existing rows 34-93 down to 36-95.
$ws.Rows("34:35").Insert()

# Populate new row 34 (Primera, 2023-02-08)
$ws.Range("A34").Value = 7
$ws.Range("B34").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C34").Value = "Ñuble"
$ws.Range("D34").Value = 44965
$ws.Range("E34").Value = 16
$ws.Range("F34").Value = 100112037
$ws.Range("G34").Value = "Cebollín"
$ws.Range("H34").Value = "Sin especificar"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 200
$ws.Range("K34").Value = 600
$ws.Range("L34").Value = 600
$ws.Range("M34").Value = 600
$ws.Range("N34").Value = "`$/paquete 6 unidades"
$ws.Range("O34").Value = "Provincia de Diguillín"
$ws.Range("P34").Value = 100
$ws.Range("Q34").Value = 6
$ws.Range("R34").Value = "Hortaliza"

# Populate new row 35 (Segunda, 2023-02-08)
$ws.Range("A35").Value = 7
$ws.Range("B35").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C35").Value = "Ñuble"
$ws.Range("D35").Value = 44965
$ws.Range("E35").Value = 16
$ws.Range("F35").Value = 100112037
$ws.Range("G35").Value = "Cebollín"
$ws.Range("H35").Value = "Sin especificar"
$ws.Range("I35").Value = "Segunda"
$ws.Range("J35").Value = 150
$ws.Range("K35").Value = 500
$ws.Range("L35").Value = 500
$ws.Range("M35").Value = 500
$ws.Range("N35").Value = "`$/paquete 6 unidades"
$ws.Range("O35").Value = "Provincia de Diguillín"
$ws.Range("P35").Value = 83
$ws.Range("Q35").Value = 6
$ws.Range("R35").Value = "Hortaliza"
